# "Add reflexion for AA" — on the "Сессии" (Sessions) sheet: fill in the
# end time of the coding session that was left open on 2023-01-04 (row 34)
# and log the start of a brand-new coding session on 2023-01-05 (row 35).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Сессии")
$ws.Activate()

# --- finish the session already in progress (row 34): fill the end time ---
$ws.Range("D34").Value = 0.84722222222222221
$ws.Range("G34").Formula = "=IF(D34<>"""",(D34-C34)*1440,0)"

# --- start a new session (row 35) ---
$ws.Range("A35").Value = "Кодирование"

$ws.Range("B35").Value = 44931
$ws.Range("B35").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("C35").Value = 0.45555555555555555
$ws.Range("C35").NumberFormat = "h:mm"

# helper formulas mirroring the pattern used by every other row
$ws.Range("F35").Formula = "=IF(I35>0,SUM(OFFSET(F36,0,0,I35,1)),0)"

$ws.Range("G35").Formula = "=IF(D35<>"""",(D35-C35)*1440,0)"
$ws.Range("G35").NumberFormat = "0"

$ws.Range("H35").Formula = "=IF(G35>0,G35-F35,0)"
$ws.Range("H35").NumberFormat = "0"

$ws.Range("I35").Formula = "=MATCH(TRUE,INDEX(((A36:A135="""")*(F36:F135="""")+(A36:A135<>"""")>0),),0)-1"

# leave the cursor where the user last clicked
$null = $ws.Range("G9").Select()
